$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "Montant" amounts from text-with-thousands-separator to real numbers
# so they can be used as source values for the pie chart.
$ws.Range("F2").Value = 84000
$ws.Range("F3").Value = 25000
$ws.Range("F4").Value = 435000
